# condition.xlsx: "Change data structure, finish shop buy/sell"
#
# Adds a new condition row (row 24) describing the "haveItem" condition
# (有可出售道具 / have sellable item), and updates the saved window/
# selection UI state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- best-effort: restore the saved window position (xWindow/yWindow in
# bookViews/workbookView). Harmless no-op if the host doesn't round-trip it.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 21320
    $win.Top = 4760
} catch {
}

# --- new row 24 data ---------------------------------------------------
# Column order: id, description, type, subtype, compareType, type2, number
# Set column B before column A so the shared-string table gets the two new
# entries in the same order as the target workbook (有可出售道具 first,
# then haveItem).
$ws.Range("B24").Value = "有可出售道具"
$ws.Range("A24").Value = "haveItem"
$ws.Range("C24").Value = "guild"
$ws.Range("D24").Value = "sellItemNumber"
$ws.Range("E24").Value = ">"
$ws.Range("F24").Value = "number"
$ws.Range("G24").Value = 0

# --- restore the saved selection ---------------------------------------
$ws.Range("A22").Select()
